$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44174
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1083
$ws.Range("D3").Value = 44169
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21000
$ws.Range("R3").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S3").Value = 1167
$ws.Range("D4").Value = 44545
$ws.Range("L4").Value = "Primera"
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 24500
$ws.Range("R4").Value = "Región de Coquimbo"
$ws.Range("S4").Value = 1361
$ws.Range("D5").Value = 44524
$ws.Range("L5").Value = "Segunda"
$ws.Range("N5").Value = 27000
$ws.Range("O5").Value = 28000
$ws.Range("P5").Value = 27500
$ws.Range("R5").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S5").Value = 1528
$ws.Range("D6").Value = 44160
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24500
$ws.Range("Q6").Value = "$/bandeja 18 kilos"
$ws.Range("R6").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S6").Value = 1361
$ws.Range("T6").Value = 18
$ws.Range("D7").Value = 44880
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 33000
$ws.Range("O7").Value = 34000
$ws.Range("P7").Value = 33500
$ws.Range("Q7").Value = "$/caja 10 kilos"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 3350
$ws.Range("T7").Value = 10
$ws.Range("D8").Value = 44544
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 21000
$ws.Range("Q8").Value = "$/bandeja 18 kilos"
$ws.Range("R8").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S8").Value = 1167
$ws.Range("T8").Value = 18
$ws.Range("D9").Value = 44881
$ws.Range("L9").Value = "Segunda"
$ws.Range("N9").Value = 41000
$ws.Range("O9").Value = 42000
$ws.Range("P9").Value = 41500
$ws.Range("R9").Value = "Región de Coquimbo"
$ws.Range("S9").Value = 2306
$ws.Range("D10").Value = 44533
$ws.Range("M10").Value = 140
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 14500
$ws.Range("Q10").Value = "$/caja 10 kilos"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 1450
$ws.Range("T10").Value = 10
